$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the existing "_GoBack" bookmark (it currently sits at the
#    very end of the document, after the last sentence of the last
#    paragraph).  It will be re-created further up, right inside the
#    "ANEXOI" -> "ANEXO I" fix below.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) Fix the "ANEXOI" typo so that it reads "ANEXO I" (i.e. insert a
#    space between the "O" and the "I").  The cursor position that
#    results from that edit is exactly where Word re-drops the
#    "_GoBack" bookmark, so we add it right after the new space.
# ------------------------------------------------------------------
$found = $d.Content.Find.Execute("XOI", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $match = $d.Content
    $match.Find.Execute("XOI", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $matchEnd = $match.End

    # Insert the missing space between "XO" and "I" (i.e. right before
    # the final "I" character of the "XOI" match).
    $rInsert = $d.Range($matchEnd - 1, $matchEnd - 1)
    $rInsert.InsertBefore(" ")

    # Re-create "_GoBack" as a zero-length bookmark right after the
    # space we just typed (between "XO " and "I"), mirroring where
    # Word leaves the insertion point after the edit.
    $rBookmark = $d.Range($matchEnd, $matchEnd)
    $d.Bookmarks.Add("_GoBack", $rBookmark)
}
